$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UCT1")
$ws.Activate()

# Row 23: change multiplier from 1.25 to 1.3 (H23, shared formula I23:AH23)
$ws.Range("H23").Formula = "=H22*1.3"
$ws.Range("I23:AH23").Formula = "=I22*1.3"

# Row 24: change multiplier from 1.8 to 2.3 (H24, shared formula I24:AH24)
$ws.Range("H24").Formula = "=H23*2.3"
$ws.Range("I24:AH24").Formula = "=I23*2.3"

# Row 25: change multiplier from 3.4 to 3.2 (H25, shared formula I25:AH25)
$ws.Range("H25").Formula = "=H21*3.2"
$ws.Range("I25:AH25").Formula = "=I21*3.2"

# Row 36: add new formula to H36
$ws.Range("H36").Formula = "=H24-H23"

# Row 37: add new formula to H37
$ws.Range("H37").Formula = "=H36/4"

# Row 38: add new formula to H38
$ws.Range("H38").Formula = "=H37+H23"

# Update sheet view: top-left cell and selection
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G37").Select()
